{"js": "// Lighting AR bug fixes: natural-language generation + grammar fixes.\n//\n// 1) \". Installation cost is estimated \" -> \". The installation cost is estimated \"\n// 2) \"...for parts and labor and the total is ${MSC}.\" ->\n//    \"...for parts and labor, resulting in a total of ${MSC}.\"\n// 3) \"The incentives are capped at 50% of the project cost and makes the\n//    modified rebate savings MRB equals to ${MRB}. Hence, the modified\n//    implementation cost (MIC) is estimated as follows:\" ->\n//    \"The incentives are capped at 50% of the project cost, which makes the\n//    modified rebate savings, MRB, equal to ${MRB}. Hence, the modified\n//    implementation cost, MIC, is estimated as follows:\"\n\nconst body = context.document.body;\n\n// --- Edit 1: \"Installation cost is estimated\" -> \"The installation cost is estimated\"\nconst r1 = body.search(\". Installation cost is estimated \", { matchCase: true });\nr1.load(\"items\");\nawait context.sync();\nif (r1.items.length !== 1) {\n  throw new Error(\"Edit 1 search expected 1 match, found \" + r1.items.length);\n}\nr1.items[0].insertText(\". The installation cost is estimated \", \"Replace\");\nawait context.sync();\n\n// --- Edit 2: motion sensor cost sentence rewording\n// NOTE: the search/replace range intentionally extends through the \"</\" that\n// follows (instead of stopping right after \"${MSC}.\") so that it fully\n// encloses the <w:proofErr w:type=\"gramStart\"/>/<w:proofErr w:type=\"gramEnd\"/>\n// pair that bracketed the old \"}.&lt;\" split. Stopping the range in the\n// middle of that marker pair (right after \"${MSC}.\") leaves a dangling,\n// unpaired gramEnd behind.\nconst r2 = body.search(\n  \"Each motion sensor costs ${MSPL} for parts and labor and the total is ${MSC}.</\",\n  { matchCase: true }\n);\nr2.load(\"items\");\nawait context.sync();\nif (r2.items.length !== 1) {\n  throw new Error(\"Edit 2 search expected 1 match, found \" + r2.items.length);\n}\nr2.items[0].insertText(\n  \"Each motion sensor costs ${MSPL} for parts and labor, resulting in a total of ${MSC}.</\",\n  \"Replace\"\n);\nawait context.sync();\n\n// --- Edit 3: incentives / rebate savings / implementation cost sentence rewording\nconst r3 = body.search(\n  \"The incentives are capped at 50% of the project cost and makes the modified rebate savings MRB equals to \",\n  { matchCase: true }\n);\nr3.load(\"items\");\nawait context.sync();\nif (r3.items.length !== 1) {\n  throw new Error(\"Edit 3a search expected 1 match, found \" + r3.items.length);\n}\nr3.items[0].insertText(\n  \"The incentives are capped at 50% of the project cost, which makes the modified rebate savings, MRB, equal to \",\n  \"Replace\"\n);\nawait context.sync();\n\nconst r4 = body.search(\n  \". Hence, the modified implementation cost (MIC) is estimated as follows:\",\n  { matchCase: true }\n);\nr4.load(\"items\");\nawait context.sync();\nif (r4.items.length !== 1) {\n  throw new Error(\"Edit 3b search expected 1 match, found \" + r4.items.length);\n}\nr4.items[0].insertText(\n  \". Hence, the modified implementation cost, MIC, is estimated as follows:\",\n  \"Replace\"\n);\nawait context.sync();\n", "ps1": "# Lighting AR bug fixes: natural-language generation + grammar fixes.\n#\n# 1) \". Installation cost is estimated \" -> \". The installation cost is estimated \"\n# 2) \"...for parts and labor and the total is ${MSC}.\" ->\n#    \"...for parts and labor, resulting in a total of ${MSC}.\"\n# 3) \"The incentives are capped at 50% of the project cost and makes the\n#    modified rebate savings MRB equals to ${MRB}. Hence, the modified\n#    implementation cost (MIC) is estimated as follows:\" ->\n#    \"The incentives are capped at 50% of the project cost, which makes the\n#    modified rebate savings, MRB, equal to ${MRB}. Hence, the modified\n#    implementation cost, MIC, is estimated as follows:\"\n\n$d = $word.ActiveDocument\n\n# wdReplace constants\n$wdReplaceNone = 0\n$wdReplaceOne = 1\n$wdReplaceAll = 2\n$wdFindStop = 0\n\n# --- Edit 1: \"Installation cost is estimated\" -> \"The installation cost is estimated\"\n$range1 = $d.Content\n$found1 = $range1.Find.Execute('. Installation cost is estimated ', $true, $false, $false, $false, $false, $true, $wdFindStop, $false, '. The installation cost is estimated ', $wdReplaceAll)\nWrite-Output \"Edit1: $found1\"\n\n# --- Edit 2: motion sensor cost sentence rewording\n# NOTE: the search/replace range intentionally extends through the \"</\" that\n# follows (instead of stopping right after \"${MSC}.\") so that it fully\n# encloses the <w:proofErr w:type=\"gramStart\"/>/<w:proofErr w:type=\"gramEnd\"/>\n# pair that bracketed the old \"}.<\" split. Stopping the range in the middle\n# of that marker pair (right after \"${MSC}.\") leaves a dangling, unpaired\n# gramEnd behind.\n$range2 = $d.Content\n$found2 = $range2.Find.Execute('Each motion sensor costs ${MSPL} for parts and labor and the total is ${MSC}.</', $true, $false, $false, $false, $false, $true, $wdFindStop, $false, 'Each motion sensor costs ${MSPL} for parts and labor, resulting in a total of ${MSC}.</', $wdReplaceAll)\nWrite-Output \"Edit2: $found2\"\n\n# --- Edit 3: incentives / rebate savings sentence rewording\n$range3 = $d.Content\n$found3 = $range3.Find.Execute('The incentives are capped at 50% of the project cost and makes the modified rebate savings MRB equals to ', $true, $false, $false, $false, $false, $true, $wdFindStop, $false, 'The incentives are capped at 50% of the project cost, which makes the modified rebate savings, MRB, equal to ', $wdReplaceAll)\nWrite-Output \"Edit3a: $found3\"\n\n# --- Edit 3b: modified implementation cost (MIC) wording\n$range4 = $d.Content\n$found4 = $range4.Find.Execute('. Hence, the modified implementation cost (MIC) is estimated as follows:', $true, $false, $false, $false, $false, $true, $wdFindStop, $false, '. Hence, the modified implementation cost, MIC, is estimated as follows:', $wdReplaceAll)\nWrite-Output \"Edit3b: $found4\"\n"}
